$d = $word.ActiveDocument
$t = $d.Tables(1)

$newValues = @(
    "37+59=96",
    "19+38=57",
    "5+36=41",
    "32+38=70",
    "40-35=5",
    "50+27=77",
    "34+57=91",
    "87-35=52",
    "87-18=69",
    "54+24=78",
    "62+36=98",
    "21+48=69",
    "66+30=96",
    "14+79=93",
    "57+27=84",
    "71-44=27",
    "81-57=24",
    "73-69=4",
    "88-39=49",
    "62-47=15",
    "55+40=95",
    "3+0=3",
    "89+4=93",
    "88+3=91",
    "85-21=64",
    "33-1=32",
    "19+36=55",
    "81+14=95",
    "85-84=1",
    "56+8=64",
    "93-71=22",
    "84-38=46",
    "74-50=24",
    "16+66=82",
    "23+4=27",
    "17-3=14",
    "84-20=64",
    "70+18=88",
    "72+7=79",
    "87-28=59",
    "14+44=58",
    "6+17=23",
    "63-37=26",
    "86-67=19",
    "29+2=31",
    "5+74=79",
    "58-4=54",
    "47-7=40",
    "72-40=32",
    "19+43=62",
    "15+9=24",
    "18+51=69",
    "32+13=45",
    "98-53=45",
    "26+70=96",
    "32+11=43",
    "42-8=34",
    "97-49=48",
    "95-50=45",
    "97-63=34",
    "12+85=97",
    "78+7=85",
    "73-21=52",
    "22-0=22",
    "90-16=74",
    "68+0=68",
    "43-40=3",
    "17+23=40",
    "20+66=86",
    "30+9=39",
    "66-36=30",
    "33+13=46",
    "38+13=51",
    "69-15=54",
    "95-32=63",
    "56-34=22",
    "86-8=78",
    "99+0=99",
    "78-16=62",
    "49+40=89",
    "25-5=20",
    "8+10=18",
    "50-26=24",
    "88-3=85",
    "75+23=98",
    "29-0=29",
    "11-5=6",
    "39+16=55",
    "5+65=70",
    "57-23=34",
    "35-1=34",
    "41-5=36",
    "23+17=40",
    "41-0=41",
    "93+5=98",
    "40+59=99",
    "34+49=83",
    "66-47=19",
    "82-79=3",
    "50-31=19"
)

$cols = $t.Columns.Count
$rows = $t.Rows.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated" $idx "cells"
